# Scheduled runner update: refresh market price snapshots (currentAveragePrice*, LevePrice*, LeveProfit*)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 388841.44
$ws.Range("I15").Value = 388841.44
$ws.Range("K15").Value = 1166524.32
$ws.Range("M15").Value = -1166355.32

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2309.7693
$ws.Range("I98").Value = 2008.45
$ws.Range("J98").Value = 3314.1667
$ws.Range("K98").Value = 2008.45
$ws.Range("L98").Value = 3314.1667
$ws.Range("M98").Value = -510.45
$ws.Range("N98").Value = -6310.1667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 98576.37
$ws.Range("I116").Value = 133982.5
$ws.Range("K116").Value = 133982.5
$ws.Range("M116").Value = -130540.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 2309.7693
$ws.Range("I122").Value = 2008.45
$ws.Range("J122").Value = 3314.1667
$ws.Range("K122").Value = 6025.35
$ws.Range("L122").Value = 9942.500100000001
$ws.Range("M122").Value = -3575.35
$ws.Range("N122").Value = -14842.5001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 5233.702
$ws.Range("I132").Value = 1810.8108
$ws.Range("J132").Value = 17898.4
$ws.Range("K132").Value = 5432.4324
$ws.Range("L132").Value = 53695.2
$ws.Range("M132").Value = -2902.4324
$ws.Range("N132").Value = -58755.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1367.5
$ws.Range("I135").Value = 538.75
$ws.Range("J135").Value = 2472.5
$ws.Range("K135").Value = 4848.75
$ws.Range("L135").Value = 22252.5
$ws.Range("M135").Value = -2313.75
$ws.Range("N135").Value = -27322.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1635.6364
$ws.Range("I2").Value = 1914
$ws.Range("J2").Value = 893.3333
$ws.Range("K2").Value = 1914
$ws.Range("L2").Value = 893.3333
$ws.Range("M2").Value = -1801
$ws.Range("N2").Value = -1119.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1635.6364
$ws.Range("I116").Value = 1914
$ws.Range("J116").Value = 893.3333
$ws.Range("K116").Value = 1914
$ws.Range("L116").Value = 893.3333
$ws.Range("M116").Value = 380
$ws.Range("N116").Value = -5481.3333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1635.6364
$ws.Range("I3").Value = 1914
$ws.Range("J3").Value = 893.3333
$ws.Range("K3").Value = 1914
$ws.Range("L3").Value = 893.3333
$ws.Range("M3").Value = -1800
$ws.Range("N3").Value = -1121.3333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1256.4
$ws.Range("I20").Value = 1222.8334
$ws.Range("J20").Value = 1287.3846
$ws.Range("K20").Value = 1222.8334
$ws.Range("L20").Value = 1287.3846
$ws.Range("M20").Value = -975.8334
$ws.Range("N20").Value = -1781.3846

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 934.6
$ws.Range("I107").Value = 885.26666
$ws.Range("J107").Value = 1082.6
$ws.Range("K107").Value = 885.26666
$ws.Range("L107").Value = 1082.6
$ws.Range("M107").Value = 1034.73334
$ws.Range("N107").Value = -4922.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1100
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 1100
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 1100
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -1674

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3141.926
$ws.Range("I31").Value = 1074.931
$ws.Range("J31").Value = 5539.64
$ws.Range("K31").Value = 1074.931
$ws.Range("L31").Value = 5539.64
$ws.Range("M31").Value = -779.931
$ws.Range("N31").Value = -6129.64

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3141.926
$ws.Range("I34").Value = 1074.931
$ws.Range("J34").Value = 5539.64
$ws.Range("K34").Value = 1074.931
$ws.Range("L34").Value = 5539.64
$ws.Range("M34").Value = -872.931
$ws.Range("N34").Value = -5943.64

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 6636.3335
$ws.Range("I94").Value = 1314.5
$ws.Range("K94").Value = 1314.5
$ws.Range("M94").Value = -863.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1100
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1100
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 1100
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -5440

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2153.8667
$ws.Range("I132").Value = 1453.5
$ws.Range("J132").Value = 3877.8462
$ws.Range("K132").Value = 4360.5
$ws.Range("L132").Value = 11633.5386
$ws.Range("M132").Value = -1830.5
$ws.Range("N132").Value = -16693.5386

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 832.6061
$ws.Range("I5").Value = 701.2143
$ws.Range("K5").Value = 2103.6429
$ws.Range("M5").Value = -1991.6429

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 586.6818
$ws.Range("I23").Value = 100
$ws.Range("J23").Value = 609.8570999999999
$ws.Range("K23").Value = 300
$ws.Range("L23").Value = 1829.5713
$ws.Range("M23").Value = -65
$ws.Range("N23").Value = -2299.5713

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 399.75
$ws.Range("J92").Value = 499.5
$ws.Range("L92").Value = 1498.5
$ws.Range("N92").Value = -3994.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 832.6061
$ws.Range("I135").Value = 701.2143
$ws.Range("K135").Value = 6310.928699999999
$ws.Range("M135").Value = -3775.928699999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6099.5103
$ws.Range("I70").Value = 3844.0557
$ws.Range("J70").Value = 12345.385
$ws.Range("K70").Value = 3844.0557
$ws.Range("L70").Value = 12345.385
$ws.Range("M70").Value = -3574.0557
$ws.Range("N70").Value = -12885.385

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 6099.5103
$ws.Range("I73").Value = 3844.0557
$ws.Range("J73").Value = 12345.385
$ws.Range("K73").Value = 3844.0557
$ws.Range("L73").Value = 12345.385
$ws.Range("M73").Value = -2908.0557
$ws.Range("N73").Value = -14217.385

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1920.7142
$ws.Range("I16").Value = 2074.1667
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 2074.1667
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -1904.1667
$ws.Range("N16").Value = -1340

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2325.0715
$ws.Range("I40").Value = 2338.7827
$ws.Range("J40").Value = 2262
$ws.Range("K40").Value = 2338.7827
$ws.Range("L40").Value = 2262
$ws.Range("M40").Value = -2202.7827
$ws.Range("N40").Value = -2534

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 215.1579
$ws.Range("I55").Value = 226.5
$ws.Range("J55").Value = 206.90909
$ws.Range("K55").Value = 226.5
$ws.Range("L55").Value = 206.90909
$ws.Range("M55").Value = -53.5
$ws.Range("N55").Value = -552.90909

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 7385.477
$ws.Range("I132").Value = 2646.6667
$ws.Range("J132").Value = 11712.218
$ws.Range("K132").Value = 7940.000100000001
$ws.Range("L132").Value = 35136.654
$ws.Range("M132").Value = -5410.000100000001
$ws.Range("N132").Value = -40196.654

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2506.5
$ws.Range("I132").Value = 1899.375
$ws.Range("J132").Value = 3316
$ws.Range("K132").Value = 5698.125
$ws.Range("L132").Value = 9948
$ws.Range("M132").Value = -3168.125
$ws.Range("N132").Value = -15008

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 21296452
$ws.Range("I136").Value = 43963052
$ws.Range("J136").Value = 443178.44
$ws.Range("K136").Value = 131889156
$ws.Range("L136").Value = 1329535.32
$ws.Range("M136").Value = -131886606
$ws.Range("N136").Value = -1334635.32
